# Atualiza os dados de teste em todas as planilhas do relatorio
# (BFS, DFS, BCU, A_Estrela_Euclidiano, A_Estrela_Haversiano):
# troca os 3 casos de teste existentes (linhas 2-4) por novos valores
# mapeados e adiciona 2 novos casos de teste (linhas 5-6).

$wb = $excel.ActiveWorkbook

$sheetNames = @("BFS", "DFS", "BCU", "A_Estrela_Euclidiano", "A_Estrela_Haversiano")

$rowsData = @{}
$rowsData["BFS"] = @(
    @{ B = 13; C = 7; D = "[13, 9, 8, 7]"; E = 16; F = 1.5; G = 0.0002608299255371094 },
    @{ B = 30; C = 23; D = "[30, 17, 20, 24, 23]"; E = 21; F = 1.238095238095238; G = 0.000125885009765625 },
    @{ B = 12; C = 26; D = "[12, 8, 5, 6, 27, 26]"; E = 29; F = 1.068965517241379; G = 0.00009179115295410156 },
    @{ B = 8; C = 18; D = "[8, 7, 10, 14, 18]"; E = 24; F = 1.041666666666667; G = 0.00008821487426757812 },
    @{ B = 23; C = 4; D = "[23, 19, 16, 12, 8, 5, 4]"; E = 25; F = 1.12; G = 0.00007867813110351562 }
)

$rowsData["DFS"] = @(
    @{ B = 13; C = 7; D = "[13, 9, 6, 3, 2, 5, 4, 7]"; E = 9; F = 0; G = 0.00007390975952148438 },
    @{ B = 30; C = 23; D = "[30, 17, 13, 9, 6, 3, 2, 5, 4, 7, 8, 12, 11, 10, 14, 15, 16, 19, 18, 21, 22, 23]"; E = 23; F = 0; G = 0.00006890296936035156 },
    @{ B = 12; C = 26; D = "[12, 8, 5, 2, 3, 6, 9, 28, 27, 26]"; E = 11; F = 0; G = 0.00006771087646484375 },
    @{ B = 8; C = 18; D = "[8, 5, 2, 3, 6, 9, 28, 29, 13, 12, 11, 10, 14, 15, 16, 17, 20, 19, 18]"; E = 25; F = 0.12; G = 0.00007486343383789062 },
    @{ B = 23; C = 4; D = "[23, 19, 16, 12, 8, 5, 4]"; E = 24; F = 0.6666666666666666; G = 0.00009202957153320312 }
)

$rowsData["BCU"] = @(
    @{ B = 13; C = 7; D = "[13, 9, 8, 7]"; E = 28; F = 0.8181818181818182; G = 0.0001370906829833984 },
    @{ B = 30; C = 23; D = "[30, 31, 32, 24, 23]"; E = 25; F = 0.7272727272727273; G = 0.0000896453857421875 },
    @{ B = 12; C = 26; D = "[12, 8, 5, 6, 27, 26]"; E = 33; F = 0.8888888888888888; G = 0.000095367431640625 },
    @{ B = 8; C = 18; D = "[8, 7, 10, 14, 18]"; E = 27; F = 0.9285714285714286; G = 0.00009441375732421875 },
    @{ B = 23; C = 4; D = "[23, 22, 21, 18, 14, 10, 7, 4]"; E = 28; F = 0.9; G = 0.0000934600830078125 }
)

$rowsData["A_Estrela_Euclidiano"] = @(
    @{ B = 13; C = 7; D = "[13, 9, 8, 7]"; E = 25; F = 3.0625; G = 0.0002100467681884766 },
    @{ B = 30; C = 23; D = "[30, 31, 32, 24, 23]"; E = 22; F = 3.0625; G = 0.0001180171966552734 },
    @{ B = 12; C = 26; D = "[12, 8, 5, 6, 27, 26]"; E = 29; F = 3.0625; G = 0.0001308917999267578 },
    @{ B = 8; C = 18; D = "[8, 7, 10, 14, 18]"; E = 24; F = 3.0625; G = 0.000118255615234375 },
    @{ B = 23; C = 4; D = "[23, 22, 21, 18, 14, 10, 7, 4]"; E = 24; F = 3.0625; G = 0.0001130104064941406 }
)

$rowsData["A_Estrela_Haversiano"] = @(
    @{ B = 13; C = 7; D = "[13, 9, 8, 7]"; E = 25; F = 3.0625; G = 0.00026702880859375 },
    @{ B = 30; C = 23; D = "[30, 31, 32, 24, 23]"; E = 22; F = 3.0625; G = 0.0001494884490966797 },
    @{ B = 12; C = 26; D = "[12, 8, 5, 6, 27, 26]"; E = 29; F = 3.0625; G = 0.0001566410064697266 },
    @{ B = 8; C = 18; D = "[8, 7, 10, 14, 18]"; E = 25; F = 3.0625; G = 0.0001356601715087891 },
    @{ B = 23; C = 4; D = "[23, 22, 21, 18, 14, 10, 7, 4]"; E = 24; F = 3.0625; G = 0.0001313686370849609 }
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $rowsData[$sheetName]

    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $i + 2
        $row = $rows[$i]

        # Linhas novas (5 e 6): preencher a coluna A ("Indice") como texto "1",
        # igual as linhas ja existentes. Linhas 2-4 ja tem esse valor, entao nao
        # precisam ser reescritas.
        if ($r -gt 4) {
            $ws.Cells.Item($r, 1).NumberFormat = "@"
            $ws.Cells.Item($r, 1).Value = "1"
        }

        $ws.Cells.Item($r, 2).Value = $row.B
        $ws.Cells.Item($r, 3).Value = $row.C
        $ws.Cells.Item($r, 4).Value = $row.D
        $ws.Cells.Item($r, 5).Value = $row.E
        $ws.Cells.Item($r, 6).Value = $row.F
        $ws.Cells.Item($r, 7).Value = $row.G
        $ws.Cells.Item($r, 8).Value = 0
    }
}
